$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1952054794520548
$ws.Range("C2").Value = 0.5547945205479452
$ws.Range("J2").Value = 0.0273972602739726
$ws.Range("P2").Value = 0.1301369863013699
$ws.Range("S2").Value = 0.09246575342465753
$ws.Range("B3").Value = 0.0223463687150838
$ws.Range("C3").Value = 0.07262569832402235
$ws.Range("J3").Value = 0.01675977653631285
$ws.Range("P3").Value = 0.7430167597765364
$ws.Range("S3").Value = 0.1452513966480447
$ws.Range("J4").Value = 0.02222222222222222
$ws.Range("P4").Value = 0.6222222222222222
$ws.Range("S4").Value = 0.3555555555555556
$ws.Range("B6").Value = 0.04958677685950413
$ws.Range("D6").Value = 0.01239669421487603
$ws.Range("F6").Value = 0.07851239669421488
$ws.Range("J6").Value = 0.2479338842975207
$ws.Range("O6").Value = 0.02892561983471074
$ws.Range("Q6").Value = 0.1735537190082645
$ws.Range("R6").Value = 0.08264462809917356
$ws.Range("S6").Value = 0.3264462809917356
$ws.Range("B7").Value = 0.08547008547008547
$ws.Range("D7").Value = 0.02136752136752137
$ws.Range("E7").Value = 0.004273504273504274
$ws.Range("F7").Value = 0.05128205128205128
$ws.Range("J7").Value = 0.1538461538461539
$ws.Range("O7").Value = 0.02564102564102564
$ws.Range("Q7").Value = 0.2008547008547009
$ws.Range("R7").Value = 0.08974358974358974
$ws.Range("S7").Value = 0.3675213675213675
$ws.Range("B8").Value = 0.07347670250896057
$ws.Range("D8").Value = 0.01254480286738351
$ws.Range("E8").Value = 0.001792114695340502
$ws.Range("J8").Value = 0.1523297491039426
$ws.Range("O8").Value = 0.02329749103942652
$ws.Range("Q8").Value = 0.1469534050179211
$ws.Range("R8").Value = 0.09139784946236559
$ws.Range("S8").Value = 0.442652329749104
$ws.Range("B9").Value = 0.09313725490196079
$ws.Range("D9").Value = 0.01470588235294118
$ws.Range("F9").Value = 0.07352941176470588
$ws.Range("J9").Value = 0.142156862745098
$ws.Range("O9").Value = 0.0196078431372549
$ws.Range("Q9").Value = 0.142156862745098
$ws.Range("R9").Value = 0.1372549019607843
$ws.Range("S9").Value = 0.3774509803921569
$ws.Range("B10").Value = 0.09457579972183588
$ws.Range("D10").Value = 0.01877607788595271
$ws.Range("E10").Value = 0.001390820584144645
$ws.Range("F10").Value = 0.06397774687065369
$ws.Range("J10").Value = 0.1300417246175243
$ws.Range("O10").Value = 0.01460361613351878
$ws.Range("Q10").Value = 0.2037552155771905
$ws.Range("R10").Value = 0.0917941585535466
$ws.Range("S10").Value = 0.3810848400556328
$ws.Range("G11").Value = 0.1391752577319588
$ws.Range("J11").Value = 0.09278350515463918
$ws.Range("K11").Value = 0.2010309278350516
$ws.Range("L11").Value = 0.5515463917525774
$ws.Range("S11").Value = 0.0154639175257732
$ws.Range("G12").Value = 0.7488584474885844
$ws.Range("J12").Value = 0.1643835616438356
$ws.Range("K12").Value = 0.0091324200913242
$ws.Range("L12").Value = 0.0365296803652968
$ws.Range("S12").Value = 0.0410958904109589
$ws.Range("G13").Value = 0.7027027027027027
$ws.Range("J13").Value = 0.2702702702702703
$ws.Range("S13").Value = 0.02702702702702703
$ws.Range("F15").Value = 0.03174603174603174
$ws.Range("H15").Value = 0.1904761904761905
$ws.Range("I15").Value = 0.05952380952380952
$ws.Range("J15").Value = 0.3373015873015873
$ws.Range("K15").Value = 0.0873015873015873
$ws.Range("M15").Value = 0.0119047619047619
$ws.Range("N15").Value = 0.003968253968253968
$ws.Range("O15").Value = 0.06349206349206349
$ws.Range("S15").Value = 0.2142857142857143
$ws.Range("F16").Value = 0.02590673575129534
$ws.Range("H16").Value = 0.1968911917098446
$ws.Range("I16").Value = 0.04145077720207254
$ws.Range("J16").Value = 0.4041450777202072
$ws.Range("K16").Value = 0.1139896373056995
$ws.Range("M16").Value = 0.03626943005181347
$ws.Range("N16").Value = 0.005181347150259068
$ws.Range("O16").Value = 0.04663212435233161
$ws.Range("S16").Value = 0.1295336787564767
$ws.Range("F17").Value = 0.01422764227642276
$ws.Range("H17").Value = 0.2052845528455285
$ws.Range("I17").Value = 0.09552845528455285
$ws.Range("J17").Value = 0.4207317073170732
$ws.Range("K17").Value = 0.08943089430894309
$ws.Range("M17").Value = 0.01829268292682927
$ws.Range("O17").Value = 0.05691056910569105
$ws.Range("S17").Value = 0.09959349593495935
$ws.Range("F18").Value = 0.01209677419354839
$ws.Range("H18").Value = 0.217741935483871
$ws.Range("I18").Value = 0.09677419354838709
$ws.Range("J18").Value = 0.3669354838709677
$ws.Range("K18").Value = 0.1008064516129032
$ws.Range("M18").Value = 0.02016129032258064
$ws.Range("N18").Value = 0.004032258064516129
$ws.Range("O18").Value = 0.0564516129032258
$ws.Range("S18").Value = 0.125
$ws.Range("F19").Value = 0.01477832512315271
$ws.Range("H19").Value = 0.227304714989444
$ws.Range("I19").Value = 0.07881773399014778
$ws.Range("J19").Value = 0.3532723434201266
$ws.Range("K19").Value = 0.131597466572836
$ws.Range("M19").Value = 0.01125967628430683
$ws.Range("N19").Value = 0.0007037297677691766
$ws.Range("O19").Value = 0.07248416608022519
$ws.Range("S19").Value = 0.1097818437719916
